$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at row 678, pushing existing rows 678-766 down to 681-769
$ws.Rows.Item(678).Insert()
$ws.Rows.Item(678).Insert()
$ws.Rows.Item(678).Insert()

# Common/static values shared by all rows in this data block
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102005
$categoria = "Naranja"
$unidad = "`$/bandeja 15 kilos granel"
$origen = "Región de O'Higgins"
$kgUnidad = 15

$fecha = Get-Date -Year 2023 -Month 8 -Day 16 -Hour 0 -Minute 0 -Second 0

# Row 678: Fukumoto / Especial
$r = 678
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Fukumoto"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 50
$ws.Cells.Item($r, 14).Value = 10000
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 10000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 667
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 679: Fukumoto / Primera
$r = 679
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Fukumoto"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 80
$ws.Cells.Item($r, 14).Value = 8000
$ws.Cells.Item($r, 15).Value = 8000
$ws.Cells.Item($r, 16).Value = 8000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 533
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 680: Fukumoto / Segunda
$r = 680
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Fukumoto"
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 100
$ws.Cells.Item($r, 14).Value = 6000
$ws.Cells.Item($r, 15).Value = 6000
$ws.Cells.Item($r, 16).Value = 6000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 400
$ws.Cells.Item($r, 20).Value = $kgUnidad

Write-Host "Final used range:" $ws.UsedRange.Address()
